# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect newly generated output data.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 36
$ws1.Range("F6").Value  = 569
$ws1.Range("F7").Value  = 1776
$ws1.Range("F8").Value  = 54
$ws1.Range("F11").Value = 2082
$ws1.Range("F12").Value = 25
$ws1.Range("F14").Value = 1337
$ws1.Range("F15").Value = 469
$ws1.Range("F17").Value = 293
$ws1.Range("F22").Value = 54
$ws1.Range("F24").Value = 12
$ws1.Range("F25").Value = 1119
$ws1.Range("F27").Value = 334
$ws1.Range("F29").Value = 271
$ws1.Range("F30").Value = 331

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 36
$ws4.Range("F6").Value  = 569
$ws4.Range("F7").Value  = 1776
$ws4.Range("F9").Value  = 54
$ws4.Range("F12").Value = 2082
$ws4.Range("F13").Value = 25
$ws4.Range("F15").Value = 1337
$ws4.Range("F16").Value = 469
$ws4.Range("F18").Value = 293
$ws4.Range("F23").Value = 54
$ws4.Range("F25").Value = 12
$ws4.Range("F26").Value = 1119
$ws4.Range("F28").Value = 334
$ws4.Range("F30").Value = 271
$ws4.Range("F31").Value = 331
